# Fix report to not trigger "emergency"
#
# Change "{%p if bad_conditions.get('emergency') %}" to
# "{%p if bad_conditions.elements.get('emergency') %}" and move the
# "_GoBack" bookmark (left behind by the last edit in Word) from its old
# location to the new edit point, right after the newly typed "elements.".

$d = $word.ActiveDocument

# Locate the text right before the insertion point ("bad_conditions.") so
# we get a stable character offset to work from.
$findRng = $d.Content
$found = $findRng.Find.Execute("{%p if bad_conditions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target template tag"
}
$boundary1 = $findRng.End

# Type the new text "elements." right after "bad_conditions."
$insPoint = $d.Range($boundary1, $boundary1)
$insPoint.InsertAfter("elements.")
$boundary2 = $boundary1 + 9   # length of "elements."

# Add a temporary bookmark at the boundary between "bad_conditions." and
# "elements." so that, once both bookmarks below exist, Word keeps the
# text split into three runs (matching how Word itself would have saved
# the document after an in-place edit followed by a bookmark move).
$tempRng = $d.Range($boundary1, $boundary1)
$d.Bookmarks.Add("_TempSplit", $tempRng)

# Re-add "_GoBack" at the new edit location (right after "elements.",
# before "get(...)"). Bookmark names are unique, so this automatically
# removes "_GoBack" from its old location near "{%p endfor %}".
$goBackRng = $d.Range($boundary2, $boundary2)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# Drop the temporary bookmark - the run split it forced stays intact.
$d.Bookmarks("_TempSplit").Delete()
